# Insert a new weekly price record for "Poroto granado" (Macroferia Regional
# de Talca) as row 94, shifting all subsequent rows down by one (dimension
# grows from A1:R131 to A1:R132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 94 - pushes rows 94..131 down to 95..132.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new record's data.
$ws.Cells.Item(94, 1).Value = 5
$ws.Cells.Item(94, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(94, 3).Value = "Maule"
$ws.Cells.Item(94, 4).Value = 44636
$ws.Cells.Item(94, 5).Value = 7
$ws.Cells.Item(94, 6).Value = 100112030
$ws.Cells.Item(94, 7).Value = "Poroto granado"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 200
$ws.Cells.Item(94, 11).Value = 18000
$ws.Cells.Item(94, 12).Value = 18000
$ws.Cells.Item(94, 13).Value = 18000
$ws.Cells.Item(94, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(94, 15).Value = "Región del Maule"
$ws.Cells.Item(94, 16).Value = 720
$ws.Cells.Item(94, 17).Value = 25
$ws.Cells.Item(94, 18).Value = "Hortaliza"
